$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell values, in natural top-to-bottom / left-to-right order so the
#     generated sharedStrings.xml table matches the original authoring
#     order (KEY, Value, UK, Birmingham.., password, Tuan@728, username,
#     mail) -----------------------------------------------------------
$ws.Range("A1").Value = "KEY"
$ws.Range("B1").Value = "Value"

$ws.Range("A2").Value = "UK"
$ws.Range("B2").Value = "Birmingham, Manchester"

$ws.Range("A3").Value = "password"
$ws.Range("B3").Value = "Tuan@728"

$ws.Range("A4").Value = "username"
$ws.Range("B4").Value = 906249919

$ws.Range("A5").Value = "mail"

# --- Formatting -----------------------------------------------------
# Header fill: accent1 theme solid fill (xl ThemeColor 5 == theme idx 4).
# Applied to A1 first, then (after the hyperlink is created) to B1 - this
# ordering reproduces the same cellXfs slot assignment as the original
# authoring session (header style = s"1", hyperlink style = s"2").
$ws.Range("A1").Interior.ThemeColor = 5

# Hyperlink on the password value cell (also auto-applies the built-in
# "Hyperlink" cell style / font to B3).
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:Tuan@728")

$ws.Range("B1").Interior.ThemeColor = 5

# Column B width (28.1 reliably round-trips to a stored OOXML width of
# exactly 29 through the engine's character -> pixel -> character
# quantization).
$ws.Columns.Item(2).ColumnWidth = 28.1

# Final selection, matching the saved cursor position in the workbook.
$ws.Range("B5").Select()
